$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23-92 down to 24-93
$ws.Rows("23").Insert()

# Populate the newly inserted row 23 with its data
$ws.Cells.Item(23, 1).Value = 11
$ws.Cells.Item(23, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(23, 3).Value = "Bíobío"
$ws.Cells.Item(23, 4).Value = 44622
$ws.Cells.Item(23, 5).Value = 8
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100101
$ws.Cells.Item(23, 8).Value = "Berries"
$ws.Cells.Item(23, 9).Value = 100101001
$ws.Cells.Item(23, 10).Value = "Arándano (blue)"
$ws.Cells.Item(23, 11).Value = "Sin especificar"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 220
$ws.Cells.Item(23, 14).Value = 2500
$ws.Cells.Item(23, 15).Value = 3000
$ws.Cells.Item(23, 16).Value = 2727
$ws.Cells.Item(23, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(23, 18).Value = "Provincia de Linares"
$ws.Cells.Item(23, 19).Value = 1364
$ws.Cells.Item(23, 20).Value = 2
